$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values are stored as text to preserve exact formatting
# (fixed decimal places, leading zeros, etc.), so force the Text number format
# before assigning, otherwise Excel would coerce the string into a float.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "236.48"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.66"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.364"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05573"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.365"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.459"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.7989"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.028"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1397"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07313"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03201"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02912"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09261"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001660"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.255"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04763"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005707"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006260"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.005067"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001052"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001499"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0004179"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.953"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.200"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04123"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007002"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.003498"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1038"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009512"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005438"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6796"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.03241"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002099"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.01009"
